$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 30 down to the new row 31 (mirrors dragging the
# fill handle down one row in the source sheet), then fill in the new data.
$ws.Range("A30:I30").Copy() | Out-Null
$ws.Range("A31:I31").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A31").Value = 45980
$ws.Range("B31").Value = 5612
$ws.Range("C31").Value = 4390
$ws.Range("D31").Value = 4065
$ws.Range("E31").Value = 236
$ws.Range("F31").Value = 54
$ws.Range("G31").Value = 26
$ws.Range("H31").Value = 9
$ws.Range("I31").Value = 0

$ws.Range("A31:I31").Select() | Out-Null
